# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.803.82"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.600.34"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.06"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -4.69%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0610"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.87"
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.822.52"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "1.605.22"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").Value = "25.804.34"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.57"
$ws.Range("D18").Value = "0.0₃0717"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.88"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.45"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.70"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.98"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.39"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "1.096.92"
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.37"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.791"
$ws.Range("E40").Value = "  -8.48%  "
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.71"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "1.735.26"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.743"
$ws.Range("E45").Value = "  -4.40%  "
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.22"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("E51").Value = "  -0.27%  "
